# WIP: update git presentation
#
# Slide 4 "Simple Git Flow": resize/reposition the screenshot picture to
# fill the space previously shared with the two command-list textboxes,
# then remove those two textboxes (id 11 "git init/add/commit" and
# id 12 "git status/log/<command> --help").
#
# Slide 5 "Branches": remove the command-list textbox (id 11, "TextBox 10":
# git branch/checkout/merge/rebase). The picture itself keeps its position.
#
# Slide 6 "Remote Repository": resize/reposition the screenshot picture,
# then remove the command-list textbox (id 11, "TextBox 10": git
# clone/fetch/push/pull/remote add origin).
#
# Slide 8 "Rebase": resize/reposition the screenshot picture, then remove
# the command-list textbox (id 2, "TextBox 1": git rebase/--onto/push
# --force).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 4 - "Simple Git Flow"
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)

$pic4 = $s4.Shapes.Item(2)
$pic4.Left = 24.215513229370117
$pic4.Top = 93.551025390625
$pic4.Width = 671.5689086914062
$pic4.Height = 255.60992431640625

# Delete the two command-list textboxes (highest index first so the
# remaining shape indices stay valid).
$s4.Shapes.Item(4).Delete()
$s4.Shapes.Item(3).Delete()

# ---------------------------------------------------------------------
# Slide 5 - "Branches"
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)

# Picture (Shapes.Item(2)) keeps its existing position/size - only the
# command-list textbox goes away.
$s5.Shapes.Item(3).Delete()

# ---------------------------------------------------------------------
# Slide 6 - "Remote Repository"
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)

$pic6 = $s6.Shapes.Item(2)
$pic6.Left = 31.529293060302734
$pic6.Top = 96.77417755126953
$pic6.Width = 658.31982421875
$pic6.Height = 276.96026611328125

$s6.Shapes.Item(3).Delete()

# ---------------------------------------------------------------------
# Slide 8 - "Rebase"
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)

$pic8 = $s8.Shapes.Item(2)
$pic8.Left = 154.0644989013672
$pic8.Top = 50.17661666870117
$pic8.Width = 527.6284790039062
$pic8.Height = 347.12396240234375

$s8.Shapes.Item(3).Delete()
